# Insert a new weekly data row at row 577 of the single data table in
# Sheet1, pushing the existing rows 577:622 down to 578:623 and filling
# the freshly inserted row with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 577 downward (this pushes old row 577 -> 578, ..., old 622 -> 623)
$ws.Rows.Item(577).Insert()

# Populate the newly inserted row 577 with the new record.
$ws.Range("A577").Value = 8
$ws.Range("B577").Value = "Terminal La Palmera de La Serena"
$ws.Range("C577").Value = "Coquimbo"
$ws.Range("D577").Value = 45223
$ws.Range("E577").Value = 4
$ws.Range("F577").Value = 100114013
$ws.Range("G577").Value = "Zanahoria"
$ws.Range("H577").Value = "Sin especificar"
$ws.Range("I577").Value = "Primera"
$ws.Range("J577").Value = 520
$ws.Range("K577").Value = 5500
$ws.Range("L577").Value = 6000
$ws.Range("M577").Value = 5750
$ws.Range("N577").Value = "`$/saco 20 kilos"
$ws.Range("O577").Value = "Provincia del Elquí"
$ws.Range("P577").Value = 288
$ws.Range("Q577").Value = 20
$ws.Range("R577").Value = "Hortaliza"
